# Drybar test data: add a new "Liquid Glass" subscribe & save product row,
# and introduce a dedicated "SubscribeSave" flag column on the DataSet sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")

# Insert a new column before AG (shifts old AG:AU -> AH:AV) and give it the
# "SubscribeSave" header that is already used as a label elsewhere in the
# workbook.
$ws.Columns("AG:AG").Insert()
$ws.Range("AG1").Value = "SubscribeSave"

# New product data row (row 62) describing the Liquid Glass subscription
# product, mirroring the style used by the other "category" rows (e.g. A58).
$ws.Range("A62").Value = "Liquid_Glass_Sub_Product"
$ws.Range("A62").Font.Color = $ws.Range("A58").Font.Color

$ws.Range("AE62").Value = "Liquid Glass High-Gloss Finishing Serum"

$ws.Range("AF62").Value = "'1"
$ws.Range("AG62").Value = "'Save 40% off daily shipments"

# Bring the DataSet sheet to the front and leave the selection on the new cell.
$ws.Activate()
$ws.Range("AA1").Select()
$ws.Range("AG62").Select()
